# Generate Report for Handoff
# Refresh the localization-status report: new handoff GUID/file names and
# updated handoff timestamps, on the Overview sheet and the per-locale
# (zh-cn / de-de) detail sheets. Hyperlink "display" text is refreshed to
# match the new file names, while the underlying hyperlink target URLs are
# preserved unchanged.

$wb = $excel.ActiveWorkbook

$oldGuidFile = "ca5c0e2d-71fc-4dfe-9259-88ff92b41384.md"
$newGuidFile = "2c179622-9275-43d8-a71d-7962e3f1b746.md"

$oldZhXlf = "ca5c0e2d-71fc-4dfe-9259-88ff92b41384.1f8b218bf63c4fd3dff78474c0739032bbb7247c.zh-cn.xlf"
$newZhXlf = "2c179622-9275-43d8-a71d-7962e3f1b746.47cd5bb7a5942e9fd083b0cd48938929924447c3.zh-cn.xlf"

$oldDeXlf = "ca5c0e2d-71fc-4dfe-9259-88ff92b41384.1f8b218bf63c4fd3dff78474c0739032bbb7247c.de-de.xlf"
$newDeXlf = "2c179622-9275-43d8-a71d-7962e3f1b746.47cd5bb7a5942e9fd083b0cd48938929924447c3.de-de.xlf"

$mdAddr = "https://github.com/OpenLocalizationTest/oltest/blob/6f79e2078409cae5e420e51bef020010d96c00d0/e2e/$oldGuidFile"
$zhXlfAddr = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a0579fdd641dadcebc3c5927570e899c9ea65849/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$oldZhXlf"
$deXlfAddr = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7c4e685bb0467902dc83290e8675d4a81f8c9f4b/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$oldDeXlf"

# ---------------------------------------------------------------------
# Overview sheet: A2 (.md link) + D2 (Latest Handoff Date)
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $mdAddr, "", "", $newGuidFile)

$wsOverview.Range("D2").Value = "2016-58-13 12:58:05"

# ---------------------------------------------------------------------
# zh-cn sheet: A2 (.md link), B2 (.md link, unchanged display),
# D2 (xlf link), E2 (Latest Handoff Datetime)
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $mdAddr, "", "", $newGuidFile)
$wsZh.Hyperlinks.Add($wsZh.Range("B2"), $mdAddr, "", "", ".md")
$wsZh.Hyperlinks.Add($wsZh.Range("D2"), $zhXlfAddr, "", "", $newZhXlf)

$wsZh.Range("E2").Value = "2016-03-13 12:58:01"

# ---------------------------------------------------------------------
# de-de sheet: A2 (.md link), B2 (.md link, unchanged display),
# D2 (xlf link), E2 (Latest Handoff Datetime)
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $mdAddr, "", "", $newGuidFile)
$wsDe.Hyperlinks.Add($wsDe.Range("B2"), $mdAddr, "", "", ".md")
$wsDe.Hyperlinks.Add($wsDe.Range("D2"), $deXlfAddr, "", "", $newDeXlf)

$wsDe.Range("E2").Value = "2016-03-13 12:58:05"
